# Eventlist.xlsx edit: add a new "bgOpcuaConnection" debugging row to the
# Table2 event list, left-align the table body, gray out the new row's
# font to flag it as newly-added/unverified, and move the selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) New data row (row 4): EventName / Data / emitted from / emitted to
$ws.Range("A4").Value = "bgOpcuaConnection"
$ws.Range("B4").Value = "1 / 0"
$ws.Range("C4").Value = "backgroundServices.js"
# D4 (emitted to) intentionally left blank.

# 2) Make the still-empty rows 4:11 inherit the same vertical/wrap format
#    that row 3 already has, so the upcoming horizontal-align change below
#    merges into the existing look (vertical-top, column B wraps + 0.00 fmt)
#    instead of creating a bare "horizontal-left-only" style.
$ws.Range("A4:D11").VerticalAlignment = -4160   # xlTop
$ws.Range("B4:B11").WrapText = $true
$ws.Range("B4:B11").NumberFormat = "0.00"

# 3) Left-align the whole remaining table body (rows 3-11).
$ws.Range("A3:D11").HorizontalAlignment = -4131 # xlLeft

# 4) Flag the freshly-added row as "not yet verified" with a muted gray font
#    (White, Background 1, Darker 35% ~ #A6A6A6 -> BGR long 10921638).
$ws.Range("A4:D4").Font.Color = 10921638

# 5) Move the selection like the author left it.
$ws.Range("A5").Select() | Out-Null
